# "If a location is used by two or more terminals, don't try to optimize it
# away" -- add a new column X that compares the "FBE" strategy results (col D)
# against the final "Add FBE strategy" results (col V), i.e. D - V, for every
# data row (3-60), the same way the existing W column ("V - B") was built:
# a plain formula on the first row, then a fill-down (shared formula) for the
# rest.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("X3").Formula = "=D3-V3"
$ws.Range("X4:X60").Formula = "=D4-V4"

# Restore the author's last selection/scroll position on the sheet.
$ws.Range("F31").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("D58").Select() | Out-Null
